$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H45").Value = 1954.1428
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 1946.6666
$ws.Range("K45").Value = 5997
$ws.Range("L45").Value = 5839.9998
$ws.Range("M45").Value = -5805
$ws.Range("N45").Value = -6223.9998
$ws.Range("H112").Value = 2796.8125
$ws.Range("J112").Value = 2950.6
$ws.Range("L112").Value = 8851.799999999999
$ws.Range("N112").Value = -11067.8
$ws.Range("H129").Value = 1203.1014
$ws.Range("J129").Value = 1352.2222
$ws.Range("L129").Value = 4056.6666
$ws.Range("N129").Value = -14056.6666
$ws.Range("H137").Value = 1492.4694
$ws.Range("I137").Value = 2054.6667
$ws.Range("J137").Value = 1070.8214
$ws.Range("K137").Value = 6164.000100000001
$ws.Range("L137").Value = 3212.4642
$ws.Range("M137").Value = -3614.000100000001
$ws.Range("N137").Value = -8312.4642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8284.102000000001
$ws.Range("I32").Value = 7303.2974
$ws.Range("K32").Value = 7303.2974
$ws.Range("M32").Value = -7016.2974
$ws.Range("H61").Value = 7241.231
$ws.Range("I61").Value = 7241.231
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7241.231
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7029.231
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 1973.6666
$ws.Range("I74").Value = 1891.1578
$ws.Range("J74").Value = 2116.182
$ws.Range("K74").Value = 1891.1578
$ws.Range("L74").Value = 2116.182
$ws.Range("M74").Value = -1017.1578
$ws.Range("N74").Value = -3864.182
$ws.Range("H76").Value = 19192
$ws.Range("J76").Value = 19192
$ws.Range("L76").Value = 19192
$ws.Range("N76").Value = -19868
$ws.Range("H77").Value = 1973.6666
$ws.Range("I77").Value = 1891.1578
$ws.Range("J77").Value = 2116.182
$ws.Range("K77").Value = 9455.789000000001
$ws.Range("L77").Value = 10580.91
$ws.Range("M77").Value = -5087.789000000001
$ws.Range("N77").Value = -19316.91
$ws.Range("H79").Value = 19192
$ws.Range("J79").Value = 19192
$ws.Range("L79").Value = 19192
$ws.Range("N79").Value = -21532
$ws.Range("H112").Value = 21001.555
$ws.Range("J112").Value = 21001.555
$ws.Range("L112").Value = 21001.555
$ws.Range("N112").Value = -23955.555
$ws.Range("H132").Value = 694047.2
$ws.Range("I132").Value = 1403287.2
$ws.Range("J132").Value = 5071.086
$ws.Range("K132").Value = 4209861.6
$ws.Range("L132").Value = 15213.258
$ws.Range("M132").Value = -4207331.6
$ws.Range("N132").Value = -20273.258
$ws.Range("H136").Value = 7241.231
$ws.Range("I136").Value = 7241.231
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 21723.693
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -19173.693
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6063157
$ws.Range("I31").Value = 2026.2069
$ws.Range("J31").Value = 12823649
$ws.Range("K31").Value = 2026.2069
$ws.Range("L31").Value = 12823649
$ws.Range("M31").Value = -1731.2069
$ws.Range("N31").Value = -12824239
$ws.Range("H34").Value = 6063157
$ws.Range("I34").Value = 2026.2069
$ws.Range("J34").Value = 12823649
$ws.Range("K34").Value = 2026.2069
$ws.Range("L34").Value = 12823649
$ws.Range("M34").Value = -1824.2069
$ws.Range("N34").Value = -12824053
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H132").Value = 3668.1667
$ws.Range("I132").Value = 2723.889
$ws.Range("J132").Value = 4612.4443
$ws.Range("K132").Value = 8171.667
$ws.Range("L132").Value = 13837.3329
$ws.Range("M132").Value = -5641.667
$ws.Range("N132").Value = -18897.3329
$ws.Range("H134").Value = 3370
$ws.Range("I134").Value = 800
$ws.Range("J134").Value = 5940
$ws.Range("K134").Value = 2400
$ws.Range("L134").Value = 17820
$ws.Range("M134").Value = 135
$ws.Range("N134").Value = -22890
$ws.Range("H140").Value = 50580
$ws.Range("J140").Value = 50580
$ws.Range("L140").Value = 50580
$ws.Range("N140").Value = -60940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 583.3333
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -3576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2500
$ws.Range("I6").Value = 2500
$ws.Range("K6").Value = 2500
$ws.Range("M6").Value = -2387
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2250
$ws.Range("H80").Value = 3188.125
$ws.Range("I80").Value = 2584.1667
$ws.Range("K80").Value = 2584.1667
$ws.Range("M80").Value = -1586.1667
$ws.Range("H83").Value = 3188.125
$ws.Range("I83").Value = 2584.1667
$ws.Range("K83").Value = 12920.8335
$ws.Range("M83").Value = -7928.833500000001
$ws.Range("H102").Value = 3573854.8
$ws.Range("I102").Value = 6495111
$ws.Range("K102").Value = 6495111
$ws.Range("M102").Value = -6493489
$ws.Range("H111").Value = 11146.5
$ws.Range("J111").Value = 11146.5
$ws.Range("L111").Value = 11146.5
$ws.Range("N111").Value = -17280.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000594.1
$ws.Range("I46").Value = 310.2
$ws.Range("J46").Value = 2000878
$ws.Range("K46").Value = 310.2
$ws.Range("L46").Value = 2000878
$ws.Range("M46").Value = -122.2
$ws.Range("N46").Value = -2001254
$ws.Range("H110").Value = 22058
$ws.Range("J110").Value = 22058
$ws.Range("L110").Value = 22058
$ws.Range("N110").Value = -30238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 713962.4399999999
$ws.Range("I136").Value = 1853084.6
$ws.Range("K136").Value = 5559253.800000001
$ws.Range("M136").Value = -5556703.800000001
